# Add a new "ViewProduct" worksheet after "Search", populate it with the
# product-lookup row, and make it the active/selected tab (mirrors the
# "test view product passed" commit).

$wb = $excel.ActiveWorkbook

$searchSheet = $wb.Worksheets.Item("Search")

# Insert the new sheet immediately after "Search" so the tab order becomes
# Authentication, Search, ViewProduct.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $searchSheet)
$ws.Name = "ViewProduct"

# Force the row-2 cells to store text (shared strings) rather than numbers,
# since several of the values ("0", "1", "2") look numeric.
$ws.Range("A2:F2").NumberFormat = "@"

$ws.Range("A2").Value = "apple watch series 7"
$ws.Range("B2").Value = "0"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "Apple Watch Series 7"

# F2 carries the same "Menlo" result-text style already used on the Search
# sheet (style index reused, no new style entries introduced).
$searchSheet.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats

# Column widths: 19.5 for column A, 25.5 for column F (the ColumnWidth
# property needs to be set ~0.8333 chars lower than the target stored
# width to account for the engine's cell-padding offset).
$ws.Columns("A").ColumnWidth = 18.666666666666668
$ws.Columns("F").ColumnWidth = 24.666666666666668

# Leave the selection on F2, matching the authored file, and make
# ViewProduct the active tab.
$ws.Range("F2").Select()
$ws.Activate()
